# Add new columns I ("I0") and J ("IF") to the stats sheet, mirroring the
# existing header style from column H, and fill in the per-row data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, border, centered/top alignment)
# from H1 into the new header cells I1:J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row (I, J) values: I = starting inning, J = finishing inning.
# For nearly every outing the pitcher started the game (I=1) and finished
# after IP innings (J = IP, the value already in column H); row 14 is a
# relief appearance that began in inning 4 and ended in inning 7.
$values = @{
    2  = @(1, 6)
    3  = @(1, 6)
    4  = @(1, 5)
    5  = @(1, 6)
    6  = @(1, 3)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 6)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 4)
    14 = @(4, 7)
    15 = @(1, 4)
    16 = @(1, 4)
    17 = @(1, 6)
    18 = @(1, 3)
    19 = @(1, 7)
    20 = @(1, 5)
    21 = @(1, 5)
    22 = @(1, 2)
}

foreach ($row in 2..22) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
